$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 224, shifting existing rows 224-237 down to 225-238
$ws.Rows.Item(224).Insert()

# Populate the newly inserted row 224 with the new weekly record
$ws.Cells.Item(224, 1).Value = 10
$ws.Cells.Item(224, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(224, 3).Value = "La Araucanía"
$ws.Cells.Item(224, 4).Value = 44585
$ws.Cells.Item(224, 5).Value = 9
$ws.Cells.Item(224, 6).Value = "Fruta"
$ws.Cells.Item(224, 7).Value = 100101
$ws.Cells.Item(224, 8).Value = "Berries"
$ws.Cells.Item(224, 9).Value = 100112025
$ws.Cells.Item(224, 10).Value = "Frutilla"
$ws.Cells.Item(224, 11).Value = "Sin especificar"
$ws.Cells.Item(224, 12).Value = "Primera"
$ws.Cells.Item(224, 13).Value = 180
$ws.Cells.Item(224, 14).Value = 7000
$ws.Cells.Item(224, 15).Value = 8000
$ws.Cells.Item(224, 16).Value = 7444
$ws.Cells.Item(224, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(224, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(224, 19).Value = 1063
$ws.Cells.Item(224, 20).Value = 7
